# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a set of rows in the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=11; I="aa"; J="Agree/Accept"},
    @{Row=13; I="aa"; J="Agree/Accept"},
    @{Row=15; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=25; I="sv"; J="Statement-opinion"},
    @{Row=45; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=47; I="sv"; J="Statement-opinion"},
    @{Row=49; I="sv"; J="Statement-opinion"},
    @{Row=57; I="sv"; J="Statement-opinion"},
    @{Row=58; I="sd"; J="Statement-non-opinion"},
    @{Row=61; I="sd"; J="Statement-non-opinion"},
    @{Row=64; I="aa"; J="Agree/Accept"},
    @{Row=66; I="sd"; J="Statement-non-opinion"},
    @{Row=69; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=70; I="sd"; J="Statement-non-opinion"},
    @{Row=73; I="qy"; J="Yes-No-Question"},
    @{Row=75; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=89; I="sv"; J="Statement-opinion"},
    @{Row=95; I="ba"; J="Appreciation"},
    @{Row=99; I="aa"; J="Agree/Accept"},
    @{Row=101; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=106; I="sv"; J="Statement-opinion"},
    @{Row=113; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=115; I="sd"; J="Statement-non-opinion"},
    @{Row=125; I="sd"; J="Statement-non-opinion"},
    @{Row=127; I="sd"; J="Statement-non-opinion"},
    @{Row=130; I="aa"; J="Agree/Accept"},
    @{Row=133; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=144; I="sv"; J="Statement-opinion"},
    @{Row=148; I="%"; J="Uninterpretable"},
    @{Row=164; I="sd"; J="Statement-non-opinion"},
    @{Row=165; I="sd"; J="Statement-non-opinion"},
    @{Row=173; I="sd"; J="Statement-non-opinion"},
    @{Row=177; I="aa"; J="Agree/Accept"},
    @{Row=178; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=198; I="sd"; J="Statement-non-opinion"},
    @{Row=207; I="sv"; J="Statement-opinion"},
    @{Row=212; I="ba"; J="Appreciation"},
    @{Row=235; I="ba"; J="Appreciation"},
    @{Row=246; I="sd"; J="Statement-non-opinion"},
    @{Row=267; I="sd"; J="Statement-non-opinion"},
    @{Row=268; I="aa"; J="Agree/Accept"},
    @{Row=277; I="aa"; J="Agree/Accept"},
    @{Row=279; I="sd"; J="Statement-non-opinion"},
    @{Row=284; I="sd"; J="Statement-non-opinion"},
    @{Row=290; I="sd"; J="Statement-non-opinion"},
    @{Row=314; I="sd"; J="Statement-non-opinion"},
    @{Row=327; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=328; I="sv"; J="Statement-opinion"},
    @{Row=333; I="aa"; J="Agree/Accept"},
    @{Row=336; I="aa"; J="Agree/Accept"},
    @{Row=354; I="sd"; J="Statement-non-opinion"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
